# Update "想去人数" (want-to-go count) figures for the latest data pull.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows, so the
# same four cells need to be bumped on each of them.

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value  = 3095
    $ws.Range("F9").Value  = 1356
    $ws.Range("F15").Value = 325
    $ws.Range("F21").Value = 2455
}
